$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 ("subtotal", with the stray pandas Series repr dumped into B14)
# is removed outright; Excel's EntireRow.Delete shifts rows 15/16 up to
# 14/15 (carrying their "recal"/"variance" labels and data with them) and
# the sheet's used range/dimension shrinks from A1:AI16 to A1:AI15.
$ws.Rows.Item(14).EntireRow.Delete()
